# Auto-generated: apply updated market/profit figures to the Sheets workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 362.63635
$ws.Range("I80").Value = 372.5
$ws.Range("J80").Value = 357
$ws.Range("K80").Value = 1117.5
$ws.Range("L80").Value = 1071
$ws.Range("M80").Value = -119.5
$ws.Range("N80").Value = -3067
$ws.Range("H83").Value = 362.63635
$ws.Range("I83").Value = 372.5
$ws.Range("J83").Value = 357
$ws.Range("K83").Value = 3352.5
$ws.Range("L83").Value = 3213
$ws.Range("M83").Value = 1639.5
$ws.Range("N83").Value = -13197
$ws.Range("H88").Value = 1950
$ws.Range("I88").Value = 900
$ws.Range("K88").Value = 900
$ws.Range("M88").Value = -494
$ws.Range("H91").Value = 1950
$ws.Range("I91").Value = 900
$ws.Range("K91").Value = 900
$ws.Range("M91").Value = 504
$ws.Range("H116").Value = 2248.75
$ws.Range("I116").Value = 1998.3334
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 1998.3334
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 1443.6666
$ws.Range("N116").Value = -9884
$ws.Range("H131").Value = 1442.6923
$ws.Range("I131").Value = 1146.25
$ws.Range("J131").Value = 5000
$ws.Range("K131").Value = 3438.75
$ws.Range("L131").Value = 15000
$ws.Range("M131").Value = 1601.25
$ws.Range("N131").Value = -25080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 773694.2
$ws.Range("I132").Value = 1047179.9
$ws.Range("J132").Value = 79461.38
$ws.Range("K132").Value = 3141539.7
$ws.Range("L132").Value = 238384.14
$ws.Range("M132").Value = -3139009.7
$ws.Range("N132").Value = -243444.14

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 224
$ws.Range("I22").Value = 224
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 224
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -51
$ws.Range("N22").ClearContents()
$ws.Range("H80").Value = 263.82352
$ws.Range("I80").Value = 224
$ws.Range("J80").Value = 280.41666
$ws.Range("K80").Value = 224
$ws.Range("L80").Value = 280.41666
$ws.Range("M80").Value = 774
$ws.Range("N80").Value = -2276.41666
$ws.Range("H83").Value = 263.82352
$ws.Range("I83").Value = 224
$ws.Range("J83").Value = 280.41666
$ws.Range("K83").Value = 1120
$ws.Range("L83").Value = 1402.0833
$ws.Range("M83").Value = 3872
$ws.Range("N83").Value = -11386.0833
$ws.Range("H86").Value = 1982
$ws.Range("I86").Value = 1448.0741
$ws.Range("J86").Value = 4041.4285
$ws.Range("K86").Value = 1448.0741
$ws.Range("L86").Value = 4041.4285
$ws.Range("M86").Value = -325.0741
$ws.Range("N86").Value = -6287.4285
$ws.Range("H89").Value = 1982
$ws.Range("I89").Value = 1448.0741
$ws.Range("J89").Value = 4041.4285
$ws.Range("K89").Value = 7240.3705
$ws.Range("L89").Value = 20207.1425
$ws.Range("M89").Value = -1624.3705
$ws.Range("N89").Value = -31439.1425

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 28000
$ws.Range("J70").Value = 28000
$ws.Range("L70").Value = 28000
$ws.Range("N70").Value = -28630
$ws.Range("H73").Value = 28000
$ws.Range("J73").Value = 28000
$ws.Range("L73").Value = 28000
$ws.Range("N73").Value = -30184
$ws.Range("H93").Value = 10922.667
$ws.Range("I93").Value = 3735.6667
$ws.Range("J93").Value = 25296.666
$ws.Range("K93").Value = 3735.6667
$ws.Range("L93").Value = 25296.666
$ws.Range("M93").Value = -1863.6667
$ws.Range("N93").Value = -29040.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 1000
$ws.Range("K13").Value = 3000
$ws.Range("M13").Value = -2832
$ws.Range("H55").Value = 980.96155
$ws.Range("J55").Value = 1245.25
$ws.Range("L55").Value = 3735.75
$ws.Range("N55").Value = -4089.75
$ws.Range("H68").Value = 731.25
$ws.Range("I68").Value = 669.2308
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 2007.6924
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -1196.6924
$ws.Range("N68").Value = -4622
$ws.Range("H71").Value = 731.25
$ws.Range("I71").Value = 669.2308
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 6023.077200000001
$ws.Range("L71").Value = 9000
$ws.Range("M71").Value = -1967.077200000001
$ws.Range("N71").Value = -17112
$ws.Range("H82").Value = 1275
$ws.Range("I82").Value = 550
$ws.Range("K82").Value = 1650
$ws.Range("M82").Value = -1244
$ws.Range("H85").Value = 1275
$ws.Range("I85").Value = 550
$ws.Range("K85").Value = 1650
$ws.Range("M85").Value = -246
$ws.Range("H122").Value = 16234666
$ws.Range("I122").Value = 19608432
$ws.Range("J122").Value = 4763865
$ws.Range("K122").Value = 176475888
$ws.Range("L122").Value = 42874785
$ws.Range("M122").Value = -176473438
$ws.Range("N122").Value = -42879685
$ws.Range("H132").Value = 32259494
$ws.Range("J132").Value = 2235.7144
$ws.Range("L132").Value = 20121.4296
$ws.Range("N132").Value = -25181.4296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 3928.8
$ws.Range("I36").Value = 675
$ws.Range("J36").Value = 5323.2856
$ws.Range("K36").Value = 675
$ws.Range("L36").Value = 5323.2856
$ws.Range("M36").Value = -190
$ws.Range("N36").Value = -6293.2856
$ws.Range("H132").Value = 1477.75
$ws.Range("I132").Value = 1187.3846
$ws.Range("J132").Value = 2017
$ws.Range("K132").Value = 3562.1538
$ws.Range("L132").Value = 6051
$ws.Range("M132").Value = -1032.1538
$ws.Range("N132").Value = -11111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1660.5
$ws.Range("I82").Value = 1356
$ws.Range("J82").Value = 2066.5
$ws.Range("K82").Value = 1356
$ws.Range("L82").Value = 2066.5
$ws.Range("M82").Value = -995
$ws.Range("N82").Value = -2788.5
$ws.Range("H85").Value = 1660.5
$ws.Range("I85").Value = 1356
$ws.Range("J85").Value = 2066.5
$ws.Range("K85").Value = 1356
$ws.Range("L85").Value = 2066.5
$ws.Range("M85").Value = -108
$ws.Range("N85").Value = -4562.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 25000
$ws.Range("J75").Value = 25000
$ws.Range("L75").Value = 25000
$ws.Range("N75").Value = -26872
$ws.Range("H78").Value = 25000
$ws.Range("J78").Value = 25000
$ws.Range("L78").Value = 75000
$ws.Range("N78").Value = -84360
$ws.Range("H81").Value = 3682.1428
$ws.Range("I81").Value = 2300
$ws.Range("J81").Value = 6170
$ws.Range("K81").Value = 4600
$ws.Range("L81").Value = 12340
$ws.Range("M81").Value = -3539
$ws.Range("N81").Value = -14462
$ws.Range("H84").Value = 3682.1428
$ws.Range("I84").Value = 2300
$ws.Range("J84").Value = 6170
$ws.Range("K84").Value = 23000
$ws.Range("L84").Value = 61700
$ws.Range("M84").Value = -17696
$ws.Range("N84").Value = -72308
$ws.Range("H132").Value = 3685.0205
$ws.Range("I132").Value = 3782.366
$ws.Range("J132").Value = 3186.125
$ws.Range("K132").Value = 11347.098
$ws.Range("L132").Value = 9558.375
$ws.Range("M132").Value = -8817.098
$ws.Range("N132").Value = -14618.375
